# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value  = 1397
$ws1.Range("F6").Value  = 68
$ws1.Range("F7").Value  = 11847
$ws1.Range("F8").Value  = 4427
$ws1.Range("F10").Value = 49
$ws1.Range("F12").Value = 19
$ws1.Range("F13").Value = 2558
$ws1.Range("F14").Value = 1106
$ws1.Range("F15").Value = 161
$ws1.Range("F16").Value = 48
$ws1.Range("F17").Value = 5137
$ws1.Range("F21").Value = 11377
$ws1.Range("F22").Value = 11348
$ws1.Range("F27").Value = 50

# --- Sheet "全部类型" ---------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value  = 1397
$ws4.Range("F6").Value  = 68
$ws4.Range("F7").Value  = 11847
$ws4.Range("F8").Value  = 4427
$ws4.Range("F10").Value = 49
$ws4.Range("F12").Value = 19
$ws4.Range("F13").Value = 2558
$ws4.Range("F15").Value = 1106
$ws4.Range("F16").Value = 161
$ws4.Range("F17").Value = 48
$ws4.Range("F18").Value = 5137
$ws4.Range("F22").Value = 11377
$ws4.Range("F23").Value = 11348
$ws4.Range("F28").Value = 50
